$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "kjhnçkinijn"
$ws.Range("F3").Select()
